$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Utilisateur inconnu" (unknown user) rows now belong to the actual user "karoui"
# (reservation created together with its ticket -- backend no longer inserts a
# placeholder "unknown user" row).
$ws.Range("A4").Value = "karoui"
$ws.Range("A5").Value = "karoui"

# New reservation row (6): same user, a fresh timestamp for the reservation
# date/time, and a "Confirmé" status highlighted like the pre-existing (but
# previously unused) green success style.
$ws.Range("A6").Value = "karoui"
$ws.Range("B6").Value = "2025-03-06T13:26:34.433205200"
$ws.Range("C6").Value = "Confirmé"

# Give the new row the same look as the rest of the table: centered dates in
# column B, reuse the existing cell formatting by copying it from a sibling
# cell so no new style entries are minted.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null

# C6 ("Confirmé") gets a green highlight. Build it on a scratch cell using the
# exact RGB of the workbook's existing (indexed) green fill so Excel reuses
# the pre-existing cellXfs/fill entry instead of minting a new one, then copy
# just the formatting onto C6.
$ws.Range("Z1").Value = "tmp"
$ws.Range("Z1").Interior.Color = 13434828
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("Z1").Clear() | Out-Null

# Column B now holds a much longer ISO timestamp string, so widen it to fit.
$ws.Columns.Item(2).ColumnWidth = 30

$excel.CutCopyMode = 0
